# Append six new weekly scoreboard rows (rows 172-177) to Sheet1, mirroring
# the formatting already used throughout the sheet (column B carries the
# short-date number format applied to every other "Date" cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A Participant, B Date, C Workout Type, D Total Duration,
# E Total Distance, F Total Elevation, G Zone1, H Zone2, I Zone3, J Zone4,
# K Zone5, L Workout Level, M Week
$newRows = @(
    @("Matt",     45480, "Walk",    174, 4.74,  1453, 150, 13,  1,  0, 0,  "Sauntering Hippo", 4),
    @("Steven",   45480, "Ride",    127, 20.08, 295,  46,  109, 11, 0, 0,  "Mighty Monkey",    4),
    @("Eric",     45480, "Workout", 91,  0,     0,    22,  48,  15, 6, 0,  "Agile Antelope",   4),
    @("Jeremiah", 45480, "Ride",    30,  11.16, 0,    0,   23,  7,  0, 0,  "Agile Antelope",   4),
    @("Jeremiah", 45480, "Workout", 24,  0,     0,    20,  4,   0,  0, 0,  "Agile Antelope",   4),
    @("Eric",     45480, "Run",     35,  3.52,  82,   0,   1,   3,  5, 25, "Agile Antelope",   4)
)

$startRow = 172
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    # Copy the number format from an existing date cell so the new date
    # values pick up the same style index the rest of column B uses,
    # instead of minting a brand-new custom numFmt.
    $ws.Range("B171").Copy() | Out-Null
    $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 2).Value = $row[1]

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
}

$excel.CutCopyMode = $false

# Move the active selection down to the new first blank row, same as Excel
# leaves the cursor after the last data entry.
$ws.Range("A178").Select() | Out-Null

"done"
